$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("54:54").Insert()
$ws.Range("A54").Borders.LineStyle = 1
$ws.Range("A54").WrapText = $true
$ws.Range("A54").VerticalAlignment = -4160
$ws.Range("A54").Font.Bold = $false
$ws.Range("A54").Font.Name = "Calibri"
$ws.Range("A54").Font.Size = 12
